# Generate Report for Handback
# Update status of d41011ff-72cf-4953-909e-a023866d6408.md (row 3 in each sheet)
# from "Ready for handoff" to "Handed back: in sync with en-US", and record
# the new handback timestamps on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is d41011ff-72cf-4953-909e-a023866d6408.md
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 is d41011ff-72cf-4953-909e-a023866d6408.md
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-22 14:51:29"

# de-de sheet: row 3 is d41011ff-72cf-4953-909e-a023866d6408.md
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-22 14:51:36"
